# Helper: write a value as literal TEXT (no numeric auto-conversion), without
# leaving a residual NumberFormat style on the cell.
function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (and
#    before "总计"), matching the fund-holdings table layout used by
#    the other quarterly sheets.
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($null, $template)
$newSheet.Name = "2022-Q1"

# NOTE: fetch the "总计" sheet only *after* the insert above, since
# inserting a sheet shifts the position of everything after it and a
# reference captured beforehand would otherwise end up stale.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy the header row formatting (bold / bordered / centered style) and
# the index-column formatting from the template sheet so the new sheet
# matches the existing look & feel.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

$template.Range("A2:A6").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)  # xlPasteFormats

# -- Header row --
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# -- Data rows --
$rows = @(
    @(0, "004995", "广发品牌消费股票A",       "2.73", "92.24", "2.94", "0.0803", 10),
    @(1, "001193", "中金消费升级股票",         "2.62", "85.49", "2.75", "0.0720", 7),
    @(2, "010245", "广发品牌消费股票C",       "0.52", "92.24", "2.94", "0.0153", 10),
    @(3, "003684", "汇安丰融灵活配置混合A", "0.01", "77.88", "3.31", "0.0003", 7),
    @(4, "003685", "汇安丰融灵活配置混合C", "0.00", "77.88", "3.31", 0,        7)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $newSheet.Cells.Item($r, 2) $row[1]
    Set-TextValue $newSheet.Cells.Item($r, 3) $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[3]
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[4]
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[5]
    if ($row[6] -eq 0) {
        $newSheet.Cells.Item($r, 7).Value = 0
    } else {
        Set-TextValue $newSheet.Cells.Item($r, 7) $row[6]
    }
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2. Insert a new summary row at the top of "总计" for the 2022-Q1
#    quarter, pushing the older rows down.
# ------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Give the new A2 index cell the same style as the other index cells
# in column A (bold / bordered / centered, style used on A3:A6).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 5
$totalSheet.Cells.Item(2,4).Value = 0.17

# Renumber the running index in column A for the rows that shifted
# down (they kept their old 0-based index values after the insert).
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(6,1).Value = 4

# Restore the original active sheet/selection (the first sheet) so the
# workbook-level view state is left as it was found.
$wb.Worksheets.Item(1).Activate()
